$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"

# ---------------------------------------------------------------------------
# 1) Split "...that the user choose to a json..." so "choose" is wrapped in
#    gramStart / gramEnd proofErr markers (grammar-checker style run split).
# ---------------------------------------------------------------------------
$pChoose = $d.Paragraphs(3)
$rChoose = $d.Range($pChoose.Range.Start, $pChoose.Range.End)
$xmlChoose = "<w:p xmlns:w='$wNs' xmlns:w14='$w14Ns' w14:paraId='7F8DEEAD' w14:textId='60347824' w:rsidR='00576A4E' w:rsidRDefault='00576A4E'>" +
    "<w:r><w:t xml:space='preserve'> - The software shall export information related to a shelter that the user </w:t></w:r>" +
    "<w:proofErr w:type='gramStart'/>" +
    "<w:r><w:t>choose</w:t></w:r>" +
    "<w:proofErr w:type='gramEnd'/>" +
    "<w:r><w:t xml:space='preserve'> to a json file when user enter shelter id and click on &#8220;Export JSON&#8221; button</w:t></w:r>" +
    "</w:p>"
$rChoose.InsertXML($xmlChoose)

# ---------------------------------------------------------------------------
# 2) Replace the empty paragraph right after "Emma" with three new
#    requirement paragraphs about xml import.
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs(8)
$rEmpty = $d.Range($pEmpty.Range.Start, $pEmpty.Range.End)
$xmlNew = "<w:p xmlns:w='$wNs'><w:r><w:t>The software must be able to accept xml files as input for putting new shelters and/or animals into the database.</w:t></w:r></w:p>" +
    "<w:p xmlns:w='$wNs'><w:r><w:t>The software should be able to recognize when the input of an xml file goes wrong and notify the user of the issue.</w:t></w:r></w:p>" +
    "<w:p xmlns:w='$wNs'><w:r><w:t>The software should prompt the user for manual input when an xml file does not include information for a shelter and/or animal that is necessary.</w:t></w:r></w:p>"
$rEmpty.InsertXML($xmlNew)

# ---------------------------------------------------------------------------
# 3) Merge the split "The software sh" + "all" + " log and save ..." runs
#    into a single run (same final text, just one <w:r>).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The software shall log and save all new inputs", $false, $false, $false, $false, $false, $true, 1, $false, "The software shall log and save all new inputs", 2)

# ---------------------------------------------------------------------------
# 4) Merge the split "The software sh" + "all" + " log all errors ..." runs
#    into a single run (same final text, just one <w:r>).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The software shall log all errors", $false, $false, $false, $false, $false, $true, 1, $false, "The software shall log all errors", 2)
